$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Relocate the "_GoBack" bookmark from its old spot (end of the
#    "(Top 5% in Computer Graphics)" bullet) to the end of the
#    "Seoul National University, 2021." bullet.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$targetRange = $d.Content
$found = $targetRange.Find.Execute("Seoul National University, 2021.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $targetRange.Collapse(0)            # wdCollapseEnd
    $targetRange.InsertAfter("X")       # temporary placeholder so the bookmark
                                         # isn't created exactly at the paragraph
                                         # boundary (avoids an engine quirk there)
    $d.Bookmarks.Add("_GoBack", $targetRange)
    $phStart = $targetRange.Start
    $placeholder = $d.Range($phStart, $phStart + 1)
    $placeholder.Delete()
}

# ------------------------------------------------------------------
# 2. Update the Advanced Optical Materials citation: add the article
#    number and change "(Accepted)" to "(Early View)".
# ------------------------------------------------------------------
$d.Content.Find.Execute(", 2025. (Accepted)", $true, $false, $false, $false, $false, $true, 1, $false, ", 2402853, 2025 (Early View). ", 2)
